$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($ws, $rangeAddr, $values) {
    $n = $values.Length
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range($rangeAddr).Value = $arr
}

# New "Total" column header (W1)
$ws.Range("W1").Value = "Total"

# Row totals (new column W) for existing category rows (2-6)
$ws.Range("W2").Value = 2155
$ws.Range("W3").Value = 288
$ws.Range("W4").Value = 729
$ws.Range("W5").Value = 395
$ws.Range("W6").Value = 1570

# New row 7: "Outros" category, with per-age-group counts and row total
$ws.Range("A7").Value = "Outros"
Set-RowValues $ws "B7:W7" @(124,5,4,44,63,102,122,195,235,282,350,454,557,531,578,493,511,413,317,97,34,5511)

# New row 8: "Total" row (column sums across all categories)
$ws.Range("A8").Value = "Total"
Set-RowValues $ws "B8:W8" @(138,8,12,53,80,123,155,254,327,424,568,767,1001,1062,1238,1126,1202,1051,728,259,72,10648)
